$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62 (ALC)
$ws.Cells.Item(62, 8).Value = 15874473
$ws.Cells.Item(62, 9).Value = 25642464
$ws.Cells.Item(62, 10).Value = 1486.5
$ws.Cells.Item(62, 11).Value = 25642464
$ws.Cells.Item(62, 12).Value = 1486.5
$ws.Cells.Item(62, 13).Value = -25641840
$ws.Cells.Item(62, 14).Value = -2734.5

# Row 64 (ALC)
$ws.Cells.Item(64, 8).Value = 4521.7393
$ws.Cells.Item(64, 9).Value = 5484.615
$ws.Cells.Item(64, 10).Value = 3270
$ws.Cells.Item(64, 11).Value = 5484.615
$ws.Cells.Item(64, 12).Value = 3270
$ws.Cells.Item(64, 13).Value = -5236.615
$ws.Cells.Item(64, 14).Value = -3766

# Row 65 (ALC)
$ws.Cells.Item(65, 8).Value = 15874473
$ws.Cells.Item(65, 9).Value = 25642464
$ws.Cells.Item(65, 10).Value = 1486.5
$ws.Cells.Item(65, 11).Value = 128212320
$ws.Cells.Item(65, 12).Value = 7432.5
$ws.Cells.Item(65, 13).Value = -128209200
$ws.Cells.Item(65, 14).Value = -13672.5

# Row 67 (ALC)
$ws.Cells.Item(67, 8).Value = 4521.7393
$ws.Cells.Item(67, 9).Value = 5484.615
$ws.Cells.Item(67, 10).Value = 3270
$ws.Cells.Item(67, 11).Value = 5484.615
$ws.Cells.Item(67, 12).Value = 3270
$ws.Cells.Item(67, 13).Value = -4626.615
$ws.Cells.Item(67, 14).Value = -4986

# Row 74 (ALC)
$ws.Cells.Item(74, 8).Value = 4543.364
$ws.Cells.Item(74, 9).Value = 4383.3335
$ws.Cells.Item(74, 10).Value = 4735.4
$ws.Cells.Item(74, 11).Value = 4383.3335
$ws.Cells.Item(74, 12).Value = 4735.4
$ws.Cells.Item(74, 13).Value = -3447.3335
$ws.Cells.Item(74, 14).Value = -6607.4

# Row 77 (ALC)
$ws.Cells.Item(77, 8).Value = 4543.364
$ws.Cells.Item(77, 9).Value = 4383.3335
$ws.Cells.Item(77, 10).Value = 4735.4
$ws.Cells.Item(77, 11).Value = 21916.6675
$ws.Cells.Item(77, 12).Value = 23677
$ws.Cells.Item(77, 13).Value = -17236.6675
$ws.Cells.Item(77, 14).Value = -33037

# Row 100 (ALC)
$ws.Cells.Item(100, 8).Value = 7248715
$ws.Cells.Item(100, 9).Value = 10418278
$ws.Cells.Item(100, 10).Value = 4000
$ws.Cells.Item(100, 11).Value = 10418278
$ws.Cells.Item(100, 12).Value = 4000
$ws.Cells.Item(100, 13).Value = -10417737
$ws.Cells.Item(100, 14).Value = -5082

# Row 137 (ALC)
$ws.Cells.Item(137, 8).Value = 1471.6552
$ws.Cells.Item(137, 9).Value = 1338.9524
$ws.Cells.Item(137, 10).Value = 1820
$ws.Cells.Item(137, 11).Value = 4016.857199999999
$ws.Cells.Item(137, 12).Value = 5460
$ws.Cells.Item(137, 13).Value = -1466.857199999999
$ws.Cells.Item(137, 14).Value = -10560

# Row 138 (ALC)
$ws.Cells.Item(138, 8).Value = 4634
$ws.Cells.Item(138, 9).Value = 1081.027
$ws.Cells.Item(138, 10).Value = 16584.908
$ws.Cells.Item(138, 11).Value = 3243.081
$ws.Cells.Item(138, 12).Value = 49754.724
$ws.Cells.Item(138, 13).Value = 1896.919
$ws.Cells.Item(138, 14).Value = -60034.724

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Cells.Item(32, 8).Value = 4324.824
$ws.Cells.Item(32, 9).Value = 3084
$ws.Cells.Item(32, 10).Value = 9026.895
$ws.Cells.Item(32, 11).Value = 3084
$ws.Cells.Item(32, 12).Value = 9026.895
$ws.Cells.Item(32, 13).Value = -2797

# Row 97 (ARM)
$ws.Cells.Item(97, 8).Value = 1322.9524
$ws.Cells.Item(97, 9).Value = 1260.6428
$ws.Cells.Item(97, 10).Value = 1447.5714
$ws.Cells.Item(97, 11).Value = 1260.6428
$ws.Cells.Item(97, 12).Value = 1447.5714
$ws.Cells.Item(97, 13).Value = -764.6428000000001
$ws.Cells.Item(97, 14).Value = -2439.5714

# Row 102 (ARM)
$ws.Cells.Item(102, 8).Value = 2059645.6
$ws.Cells.Item(102, 9).Value = 2471174.8
$ws.Cells.Item(102, 10).Value = 2000
$ws.Cells.Item(102, 11).Value = 2471174.8
$ws.Cells.Item(102, 12).Value = 2000
$ws.Cells.Item(102, 13).Value = -2469552.8

# Row 132 (ARM)
$ws.Cells.Item(132, 8).Value = 2231.0576
$ws.Cells.Item(132, 9).Value = 1404.7805
$ws.Cells.Item(132, 10).Value = 5310.8184
$ws.Cells.Item(132, 11).Value = 4214.3415
$ws.Cells.Item(132, 12).Value = 15932.4552
$ws.Cells.Item(132, 13).Value = -1684.3415
$ws.Cells.Item(132, 14).Value = -20992.4552

$ws = $wb.Worksheets.Item("BSM")
# Row 44 (BSM)
$ws.Cells.Item(44, 8).Value = 18100
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 18100
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 12).Value = 18100
$ws.Cells.Item(44, 14).Value = -19094

# Row 86 (BSM)
$ws.Cells.Item(86, 8).Value = 27779412
$ws.Cells.Item(86, 9).Value = 33334884
$ws.Cells.Item(86, 10).Value = 2050
$ws.Cells.Item(86, 11).Value = 33334884
$ws.Cells.Item(86, 12).Value = 2050
$ws.Cells.Item(86, 13).Value = -33333761

# Row 89 (BSM)
$ws.Cells.Item(89, 8).Value = 27779412
$ws.Cells.Item(89, 9).Value = 33334884
$ws.Cells.Item(89, 10).Value = 2050
$ws.Cells.Item(89, 11).Value = 166674420
$ws.Cells.Item(89, 12).Value = 10250
$ws.Cells.Item(89, 13).Value = -166668804

# Row 94 (BSM)
$ws.Cells.Item(94, 8).Value = 2307.0715
$ws.Cells.Item(94, 9).Value = 1974.875
$ws.Cells.Item(94, 10).Value = 2750
$ws.Cells.Item(94, 11).Value = 1974.875
$ws.Cells.Item(94, 12).Value = 2750
$ws.Cells.Item(94, 13).Value = -1523.875
$ws.Cells.Item(94, 14).Value = -3652

# Row 99 (BSM)
$ws.Cells.Item(99, 8).Value = 76924260
$ws.Cells.Item(99, 9).Value = 125001130
$ws.Cells.Item(99, 10).Value = 1260
$ws.Cells.Item(99, 11).Value = 125001130
$ws.Cells.Item(99, 12).Value = 1260
$ws.Cells.Item(99, 13).Value = -124999632
$ws.Cells.Item(99, 14).Value = -4256

# Row 105 (BSM)
$ws.Cells.Item(105, 8).Value = 1930
$ws.Cells.Item(105, 9).Value = 1671.4286
$ws.Cells.Item(105, 10).Value = 2533.3333
$ws.Cells.Item(105, 11).Value = 1671.4286
$ws.Cells.Item(105, 12).Value = 2533.3333
$ws.Cells.Item(105, 13).Value = 75.57140000000004
$ws.Cells.Item(105, 14).Value = -6027.3333

$ws = $wb.Worksheets.Item("CRP")
# Row 45 (CRP)
$ws.Cells.Item(45, 8).Value = 7500
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 10).Value = 7500
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 12).Value = 7500
$ws.Cells.Item(45, 14).Value = -8686

# Row 58 (CRP)
$ws.Cells.Item(58, 8).Value = 1604.6666
$ws.Cells.Item(58, 9).Value = 990.2857
$ws.Cells.Item(58, 10).Value = 2321.4443
$ws.Cells.Item(58, 11).Value = 990.2857
$ws.Cells.Item(58, 12).Value = 2321.4443
$ws.Cells.Item(58, 13).Value = -787.2857
$ws.Cells.Item(58, 14).Value = -2727.4443

# Row 62 (CRP)
$ws.Cells.Item(62, 8).Value = 9035
$ws.Cells.Item(62, 9).Value = 9035
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 9035
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).Value = -8411
$ws.Cells.Item(62, 14).ClearContents()

# Row 65 (CRP)
$ws.Cells.Item(65, 8).Value = 9035
$ws.Cells.Item(65, 9).Value = 9035
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 45175
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = -42055
$ws.Cells.Item(65, 14).ClearContents()

# Row 105 (CRP)
$ws.Cells.Item(105, 8).Value = 19610348
$ws.Cells.Item(105, 9).Value = 27780534
$ws.Cells.Item(105, 10).Value = 1900
$ws.Cells.Item(105, 11).Value = 27780534
$ws.Cells.Item(105, 12).Value = 1900
$ws.Cells.Item(105, 13).Value = -27778787
$ws.Cells.Item(105, 14).Value = -5394

# Row 122 (CRP)
$ws.Cells.Item(122, 8).Value = 2601.1
$ws.Cells.Item(122, 9).Value = 1867.7142
$ws.Cells.Item(122, 10).Value = 4312.3335
$ws.Cells.Item(122, 11).Value = 5603.142599999999
$ws.Cells.Item(122, 12).Value = 12937.0005
$ws.Cells.Item(122, 13).Value = -3153.142599999999
$ws.Cells.Item(122, 14).Value = -17837.0005

# Row 134 (CRP)
$ws.Cells.Item(134, 8).Value = 2041.28
$ws.Cells.Item(134, 9).Value = 2197.081
$ws.Cells.Item(134, 10).Value = 1597.8462
$ws.Cells.Item(134, 11).Value = 6591.243
$ws.Cells.Item(134, 12).Value = 4793.5386
$ws.Cells.Item(134, 13).Value = -4056.243
$ws.Cells.Item(134, 14).Value = -9863.5386

# Row 136 (CRP)
$ws.Cells.Item(136, 8).Value = 1604.6666
$ws.Cells.Item(136, 9).Value = 990.2857
$ws.Cells.Item(136, 10).Value = 2321.4443
$ws.Cells.Item(136, 11).Value = 2970.8571
$ws.Cells.Item(136, 12).Value = 6964.3329
$ws.Cells.Item(136, 13).Value = -420.8571000000002
$ws.Cells.Item(136, 14).Value = -12064.3329

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (GSM)
$ws.Cells.Item(70, 8).Value = 4714.9
$ws.Cells.Item(70, 9).Value = 4006.4443
$ws.Cells.Item(70, 10).Value = 5294.5454
$ws.Cells.Item(70, 11).Value = 4006.4443
$ws.Cells.Item(70, 12).Value = 5294.5454
$ws.Cells.Item(70, 13).Value = -3736.4443
$ws.Cells.Item(70, 14).Value = -5834.5454

# Row 73 (GSM)
$ws.Cells.Item(73, 8).Value = 4714.9
$ws.Cells.Item(73, 9).Value = 4006.4443
$ws.Cells.Item(73, 10).Value = 5294.5454
$ws.Cells.Item(73, 11).Value = 4006.4443
$ws.Cells.Item(73, 12).Value = 5294.5454
$ws.Cells.Item(73, 13).Value = -3070.4443
$ws.Cells.Item(73, 14).Value = -7166.5454

# Row 80 (GSM)
$ws.Cells.Item(80, 8).Value = 2521.5625
$ws.Cells.Item(80, 9).Value = 2361.4285
$ws.Cells.Item(80, 10).Value = 2827.2727
$ws.Cells.Item(80, 11).Value = 2361.4285
$ws.Cells.Item(80, 12).Value = 2827.2727
$ws.Cells.Item(80, 13).Value = -1363.4285
$ws.Cells.Item(80, 14).Value = -4823.2727

# Row 83 (GSM)
$ws.Cells.Item(83, 8).Value = 2521.5625
$ws.Cells.Item(83, 9).Value = 2361.4285
$ws.Cells.Item(83, 10).Value = 2827.2727
$ws.Cells.Item(83, 11).Value = 11807.1425
$ws.Cells.Item(83, 12).Value = 14136.3635
$ws.Cells.Item(83, 13).Value = -6815.1425
$ws.Cells.Item(83, 14).Value = -24120.3635

# Row 97 (GSM)
$ws.Cells.Item(97, 8).Value = 802.5
$ws.Cells.Item(97, 9).Value = 836.6667
$ws.Cells.Item(97, 10).Value = 700
$ws.Cells.Item(97, 11).Value = 836.6667
$ws.Cells.Item(97, 12).Value = 700
$ws.Cells.Item(97, 13).Value = -340.6667

# Row 132 (GSM)
$ws.Cells.Item(132, 8).Value = 2008.6757
$ws.Cells.Item(132, 9).Value = 1651.4333
$ws.Cells.Item(132, 10).Value = 3539.7144
$ws.Cells.Item(132, 11).Value = 4954.2999
$ws.Cells.Item(132, 12).Value = 10619.1432
$ws.Cells.Item(132, 13).Value = -2424.2999

$ws = $wb.Worksheets.Item("LTW")
# Row 93 (LTW)
$ws.Cells.Item(93, 8).Value = 21740168
$ws.Cells.Item(93, 9).Value = 948.8333
$ws.Cells.Item(93, 10).Value = 100001360
$ws.Cells.Item(93, 11).Value = 948.8333
$ws.Cells.Item(93, 12).Value = 100001360
$ws.Cells.Item(93, 13).Value = 299.1667
$ws.Cells.Item(93, 14).Value = -100003856

# Row 100 (LTW)
$ws.Cells.Item(100, 8).Value = 1522.2
$ws.Cells.Item(100, 9).Value = 1465.25
$ws.Cells.Item(100, 10).Value = 1750
$ws.Cells.Item(100, 11).Value = 1465.25
$ws.Cells.Item(100, 12).Value = 1750
$ws.Cells.Item(100, 13).Value = -924.25

# Row 136 (LTW)
$ws.Cells.Item(136, 8).Value = 7463.268
$ws.Cells.Item(136, 9).Value = 5224.032
$ws.Cells.Item(136, 10).Value = 14404.9
$ws.Cells.Item(136, 11).Value = 15672.096
$ws.Cells.Item(136, 12).Value = 43214.7
$ws.Cells.Item(136, 13).Value = -13122.096
$ws.Cells.Item(136, 14).Value = -48314.7

$ws = $wb.Worksheets.Item("WVR")
# Row 62 (WVR)
$ws.Cells.Item(62, 8).Value = 2900
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 2900
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 2900
$ws.Cells.Item(62, 13).ClearContents()
$ws.Cells.Item(62, 14).Value = -4148

# Row 65 (WVR)
$ws.Cells.Item(65, 8).Value = 2900
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 2900
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 14500
$ws.Cells.Item(65, 13).ClearContents()
$ws.Cells.Item(65, 14).Value = -20740

# Row 81 (WVR)
$ws.Cells.Item(81, 8).Value = 1720.2
$ws.Cells.Item(81, 9).Value = 900.25
$ws.Cells.Item(81, 10).Value = 5000
$ws.Cells.Item(81, 11).Value = 1800.5
$ws.Cells.Item(81, 12).Value = 10000
$ws.Cells.Item(81, 13).Value = -739.5
$ws.Cells.Item(81, 14).Value = -12122

# Row 84 (WVR)
$ws.Cells.Item(84, 8).Value = 1720.2
$ws.Cells.Item(84, 9).Value = 900.25
$ws.Cells.Item(84, 10).Value = 5000
$ws.Cells.Item(84, 11).Value = 9002.5
$ws.Cells.Item(84, 12).Value = 50000
$ws.Cells.Item(84, 13).Value = -3698.5
$ws.Cells.Item(84, 14).Value = -60608

# Row 132 (WVR)
$ws.Cells.Item(132, 8).Value = 989.6353
$ws.Cells.Item(132, 9).Value = 651.34424
$ws.Cells.Item(132, 10).Value = 1849.4584
$ws.Cells.Item(132, 11).Value = 1954.03272
$ws.Cells.Item(132, 12).Value = 5548.3752
$ws.Cells.Item(132, 13).Value = 575.9672799999998
$ws.Cells.Item(132, 14).Value = -10608.3752

Write-Host "Applied 41 row updates"